# Refresh the crypto price list ("Updated cryptos list ... with GitHub Actions").
# Column D (Price) and E (Volume(1h)) get refreshed numbers for most coins;
# two pairs of rows (38/39 "Maker" <-> "MXToken" and 50/51 "Aptos"/"Cronos" ->
# "Cronos"/"Algorand") had their whole row's data re-ranked/replaced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a cell as literal text, never letting Excel's automatic
# number/date recognition mangle values like "1.000", "5.910" or
# "0.00000000118" (which would otherwise collapse to 1, 5.91, 1.18E-09).
# A leading apostrophe forces text entry; resetting the style back to
# "Normal" afterwards drops the quote-prefix formatting Excel applies so the
# cell's style stays identical to how it started.
function Set-TextValue([string]$addr, [string]$value) {
    $rng = $ws.Range($addr)
    $rng.Value = "'" + $value
    $rng.Style = "Normal"
}

# row 2 - Bitcoin
Set-TextValue "D2" "29.230.75"
$ws.Range("E2").Value = "  -0.42%  "

# row 3 - Ethereum
Set-TextValue "D3" "1.830.51"

# row 4 - TetherUSD
Set-TextValue "D4" "1.000"
$ws.Range("E4").Value = "  +0.16%  "

# row 5 - BNB
Set-TextValue "D5" "236.63"
$ws.Range("E5").Value = "  -1.43%  "

# row 6 - XRP
Set-TextValue "D6" "0.6069"
$ws.Range("E6").Value = "  -3.70%  "

# row 7 - USDC
$ws.Range("E7").Value = "  +0.16%  "

# row 8 - Dogecoin
Set-TextValue "D8" "0.07113"
$ws.Range("E8").Value = "  -4.83%  "

# row 9 - Cardano
Set-TextValue "D9" "0.2819"
$ws.Range("E9").Value = "  -2.85%  "

# row 10 - Solana
Set-TextValue "D10" "23.95"
$ws.Range("E10").Value = "  -3.98%  "

# row 11 - TRON
Set-TextValue "D11" "0.07671"
$ws.Range("E11").Value = "  -0.72%  "

# row 12 - WrappedEther
Set-TextValue "D12" "1.830.85"
$ws.Range("E12").Value = "  -0.72%  "

# row 13 - Polkadot
Set-TextValue "D13" "4.833"
$ws.Range("E13").Value = "  -2.94%  "

# row 14 - ShibaInu
$ws.Range("E14").Value = "  -1.94%  "

# row 15 - Polygon
Set-TextValue "D15" "0.6376"
$ws.Range("E15").Value = "  -5.93%  "

# row 16 - WrappedliquidstakedEther2.0
Set-TextValue "D16" "2.079.44"
$ws.Range("E16").Value = "  -0.60%  "

# row 17 - Litecoin
Set-TextValue "D17" "79.52"
$ws.Range("E17").Value = "  -2.95%  "

# row 18 - Uniswap
Set-TextValue "D18" "5.910"
$ws.Range("E18").Value = "  -5.11%  "

# row 19 - WrappedBTC
Set-TextValue "D19" "29.218.46"
$ws.Range("E19").Value = "  -0.45%  "

# row 20 - BitcoinCash
Set-TextValue "D20" "228.99"
$ws.Range("E20").Value = "  -0.14%  "

# row 21 - Avalanche
$ws.Range("E21").Value = "  -3.91%  "

# row 22 - Dai
$ws.Range("E22").Value = "  +0.15%  "

# row 23 - Chainlink
Set-TextValue "D23" "7.033"
$ws.Range("E23").Value = "  -4.84%  "

# row 24 - BinanceUSD
Set-TextValue "D24" "1.001"
$ws.Range("E24").Value = "  +0.05%  "

# row 25 - Monero
Set-TextValue "D25" "154.41"
$ws.Range("E25").Value = "  -2.34%  "

# row 26 - Cosmos
Set-TextValue "D26" "8.091"
$ws.Range("E26").Value = "  -5.19%  "

# row 27 - Stellar
Set-TextValue "D27" "0.1294"
$ws.Range("E27").Value = "  -4.59%  "

# row 28 - EthereumClassic
Set-TextValue "D28" "16.66"
$ws.Range("E28").Value = "  -4.67%  "

# row 29 - Toncoin
Set-TextValue "D29" "1.489"
$ws.Range("E29").Value = "  +2.20%  "

# row 30 - Hedera
Set-TextValue "D30" "0.06479"
$ws.Range("E30").Value = "  -5.72%  "

# row 31 - PancakeSwap
Set-TextValue "D31" "1.459"
$ws.Range("E31").Value = "  -2.08%  "

# row 32 - Filecoin
Set-TextValue "D32" "3.832"
$ws.Range("E32").Value = "  -5.72%  "

# row 33 - InternetComputer(DFINITY)
Set-TextValue "D33" "3.830"
$ws.Range("E33").Value = "  -5.95%  "

# row 34 - ARBITRUM
$ws.Range("E34").Value = "  -0.96%  "

# row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -4.72%  "

# row 36 - ImmutableX
Set-TextValue "D36" "0.6538"
$ws.Range("E36").Value = "  -6.67%  "

# row 37 - HuobiToken
Set-TextValue "D37" "2.558"
$ws.Range("E37").Value = "  -1.02%  "

# row 38 - was "Maker", now "MXToken" (rows 38/39 swapped rank)
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D38" "2.759"
$ws.Range("E38").Value = "  -2.19%  "

# row 39 - was "MXToken", now "Maker"
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D39" "1.221.76"
$ws.Range("E39").Value = "  -1.28%  "

# row 40 - VeChain
Set-TextValue "D40" "0.01750"
$ws.Range("E40").Value = "  -5.21%  "

# row 41 - FraxShare
Set-TextValue "D41" "6.517"
$ws.Range("E41").Value = "  -4.05%  "

# row 42 - TrustWalletToken
Set-TextValue "D42" "0.9321"
$ws.Range("E42").Value = "  -1.32%  "

# row 44 - Quant
Set-TextValue "D44" "101.14"
$ws.Range("E44").Value = "  +0.16%  "

# row 45 - RocketPoolETH
Set-TextValue "D45" "1.979.17"
$ws.Range("E45").Value = "  -1.14%  "

# row 46 - Aave
Set-TextValue "D46" "63.45"
$ws.Range("E46").Value = "  -3.12%  "

# row 47 - BabyDogeCoin
Set-TextValue "D47" "0.00000000118"
$ws.Range("E47").Value = "  -1.80%  "

# row 48 - RenderToken
Set-TextValue "D48" "1.612"
$ws.Range("E48").Value = "  -6.28%  "

# row 49 - EnergySwap
Set-TextValue "D49" "8.546"
$ws.Range("E49").Value = "  -4.78%  "

# row 50 - was "Aptos", now "Cronos"
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D50" "0.05543"
$ws.Range("E50").Value = "  -2.34%  "

# row 51 - was "Cronos", now "Algorand"
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D51" "0.1078"
$ws.Range("E51").Value = "  -5.87%  "
